$d = $word.ActiveDocument

# 1. "As a Data Scientist at Mintek, I " -> "As a Data Scientist intern at Mintek, I "
$d.Content.Find.Execute("As a Data Scientist at Mintek, I ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "As a Data Scientist intern at Mintek, I ", 2)

# 2. Expand the KPI dashboard sentence
$old2 = "Additionally, I design and maintain organisational Key Performance Indicators (KPIs) dashboards to support real-time tracking and performance monitoring. I also execute"
$new2 = "Additionally, I design and maintain organisational Key Performance Indicators (KPIs) dashboards, including research outputs, HR metrics, divisional performance, and other key areas to support real-time tracking and performance monitoring, reducing manual reporting time by 40%. I also execute"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new2, 2)

# 3. Remove lastRenderedPageBreak before "For NEXT TIME" - handled automatically by engine normally.

# 4. Machine Learning / PyTorch - proofErr tags removal only; text unchanged.

Write-Output "done"
